$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the B2:D9 block to 0 by default
$ws.Range("B2:D9").Value = 0

# Apply the specific non-zero exceptions per the diff
$ws.Range("C3").Value = -0.6407510497588037
$ws.Range("C7").Value = -0.6788932807695982
